# Weekly update: insert two new "Fruta" / "Ciruela" price rows (variety
# "Fortuna", sampled 2023-04-05 / serial 45021) at the top of the existing
# Terminal Hortofrutícola Agro Chillán block (rows 106-115), pushing the
# pre-existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 106 - everything that used to
# live in rows 106:115 shifts down to 108:117 automatically, carrying all of
# its cell content and formatting with it.
$ws.Range("A106:A107").EntireRow.Insert()

# ---- New row 106: Fortuna / Especial --------------------------------------
$ws.Cells.Item(106, 1).Value2 = 7
$ws.Cells.Item(106, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(106, 3).Value2 = "Ñuble"
$ws.Cells.Item(106, 4).Value2 = 45021
$ws.Cells.Item(106, 5).Value2 = 16
$ws.Cells.Item(106, 6).Value2 = "Fruta"
$ws.Cells.Item(106, 7).Value2 = 100103
$ws.Cells.Item(106, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(106, 9).Value2 = 100103002
$ws.Cells.Item(106, 10).Value2 = "Ciruela"
$ws.Cells.Item(106, 11).Value2 = "Fortuna"
$ws.Cells.Item(106, 12).Value2 = "Especial"
$ws.Cells.Item(106, 13).Value2 = 60
$ws.Cells.Item(106, 14).Value2 = 12000
$ws.Cells.Item(106, 15).Value2 = 12000
$ws.Cells.Item(106, 16).Value2 = 12000
$ws.Cells.Item(106, 17).Value2 = "$/bandeja 18 kilos granel"
$ws.Cells.Item(106, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(106, 19).Value2 = 667
$ws.Cells.Item(106, 20).Value2 = 18

# ---- New row 107: Fortuna / Primera ----------------------------------------
$ws.Cells.Item(107, 1).Value2 = 7
$ws.Cells.Item(107, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(107, 3).Value2 = "Ñuble"
$ws.Cells.Item(107, 4).Value2 = 45021
$ws.Cells.Item(107, 5).Value2 = 16
$ws.Cells.Item(107, 6).Value2 = "Fruta"
$ws.Cells.Item(107, 7).Value2 = 100103
$ws.Cells.Item(107, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(107, 9).Value2 = 100103002
$ws.Cells.Item(107, 10).Value2 = "Ciruela"
$ws.Cells.Item(107, 11).Value2 = "Fortuna"
$ws.Cells.Item(107, 12).Value2 = "Primera"
$ws.Cells.Item(107, 13).Value2 = 60
$ws.Cells.Item(107, 14).Value2 = 10000
$ws.Cells.Item(107, 15).Value2 = 10000
$ws.Cells.Item(107, 16).Value2 = 10000
$ws.Cells.Item(107, 17).Value2 = "$/bandeja 18 kilos granel"
$ws.Cells.Item(107, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(107, 19).Value2 = 556
$ws.Cells.Item(107, 20).Value2 = 18

"OK"
